# Agrupamento de distribuidoras ENF+EMG=EMR e EBO+EPB=EPB
# Adds the new aggregated distributor "EMR" (MG / SE / SE) as a new
# row at the bottom of the "Plan1" lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row right after the current last data row (row 56 -> 57)
$ws.Range("A57").Value = "EMR"
$ws.Range("B57").Value = "MG"
$ws.Range("C57").Value = "SE"
$ws.Range("D57").Value = "SE"

# Match the reviewer's scroll position / selection after the edit
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A58").Select()
